# V 2.0.2 se arreglo la fecha y hora de reimpresion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Paciente: Apellidos / Nombres / No. Expediente Clinico
$ws.Range("A6").Value = "MENDEZ"
$ws.Range("C6").Value = "POP"
$ws.Range("E6").Value = "HANSEL MARIO"
$ws.Range("G6").Value = "ALEXANDER"
$ws.Range("I6").Value = "2017-26249/201762628"

# Direccion actual (Calle, Municipio, Telefono)
$ws.Range("A8").Value = "BARRIO SAN JUAN"
$ws.Range("D8").Value = "SAN AGUSTIN ACASAGUASTLAN EL PROGRESO"
$ws.Range("J8").Value = "31351928"

# Direccion habitual
$ws.Range("A10").Value = ""

# Fecha de nacimiento / Edad / Lugar de nacimiento
$ws.Range("A12").Value = "1998-04-09"
$ws.Range("F12").Value = "19"
$ws.Range("H12").Value = "JALAPA"

# Estado civil / Ocupacion / Nacionalidad / No. de Cedula
$ws.Range("A14").Value = "Soltero"
$ws.Range("D14").Value = "ESTUDIANTE"
$ws.Range("F14").Value = "GUATEMALTECO"
$ws.Range("H14").Value = ""

# Nombre del Conyugue
$ws.Range("A16").Value = ""

# Nombre del Padre / Nombre de la Madre
$ws.Range("A18").Value = "LUIS MENDEZ"
$ws.Range("F18").Value = "GLENIA POP"

# En caso de emergencia: Nombre / Parentesco / Direccion / Telefono
$ws.Range("A20").Value = "GLENIA POP"
$ws.Range("F20").Value = "MADRE"
$ws.Range("H20").Value = ""
$ws.Range("J20").Value = "31351928"

# Fecha de Ingreso / Hora / Servicio
$ws.Range("A24").Value = "24/10/2017"
$ws.Range("C24").Value = "14:51:18"
$ws.Range("D24").Value = "HEMATO-ONCO"
